# Append the latest EUR -> ARS quote as a new row at the bottom of the
# sheet (row 57), extending the used range from A1:C56 to A1:C57.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

# Column A holds a date-like string ("2025-10-04"). Force text formatting
# before assigning so Excel stores it as a literal string (matching every
# other row in the sheet) instead of auto-converting it to a date serial.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-04"

$ws.Cells.Item($row, 2).Value = "15:18:28"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,794.1737"
